$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.914
$ws.Range("C3").Value = -12.244
$ws.Range("E8").Value = 16.715
$ws.Range("E11").Value = 16.789
$ws.Range("A12").Value = -21.654
$ws.Range("B14").Value = 5.828
$ws.Range("E14").Value = 16.977
$ws.Range("E15").Value = 16.177
$ws.Range("E17").Value = 16.647
$ws.Range("C20").Value = -12.46
$ws.Range("C25").Value = -12.738
$ws.Range("B26").Value = 6.548
$ws.Range("E26").Value = 16.506
$ws.Range("A27").Value = -21.531
$ws.Range("C30").Value = -13.274
$ws.Range("B31").Value = 6.167999999999999
$ws.Range("A32").Value = -21.495
$ws.Range("B35").Value = 7.342000000000001
$ws.Range("A36").Value = -21.333
$ws.Range("E36").Value = 16.868
$ws.Range("B37").Value = 7.641
$ws.Range("A38").Value = -20.363
$ws.Range("C44").Value = -12.199
$ws.Range("B45").Value = 5.747
$ws.Range("A46").Value = -21.591
$ws.Range("C47").Value = -12.321
$ws.Range("B52").Value = 4.743
$ws.Range("A54").Value = -21.585
$ws.Range("A55").Value = -21.825
$ws.Range("A56").Value = -21.768
$ws.Range("B57").Value = 6.090000000000001
$ws.Range("C58").Value = -13.289
$ws.Range("E64").Value = 17.254
$ws.Range("A67").Value = -21.6
$ws.Range("A69").Value = -21.651
$ws.Range("A72").Value = -21.567
$ws.Range("C78").Value = -13.032
$ws.Range("E79").Value = 17.398
$ws.Range("B81").Value = 6.693
$ws.Range("A83").Value = -20.489
$ws.Range("B83").Value = 7.256
$ws.Range("C84").Value = -13.509
$ws.Range("A86").Value = -22.292
$ws.Range("C89").Value = -11.045
$ws.Range("E89").Value = 17.163
$ws.Range("A91").Value = -21.768
$ws.Range("C91").Value = -11.156
$ws.Range("C92").Value = -11.166
$ws.Range("A93").Value = -21.593
$ws.Range("C96").Value = -13.109
$ws.Range("A99").Value = -20.683
$ws.Range("B100").Value = 5.765
$ws.Range("B102").Value = 6.667
$ws.Range("C102").Value = -12.798
